# "updated new logo, now use new logo in index"
#
# Re-positions/re-sizes the "Rectangle 5" logo-text shape on slide 2 and
# shrinks its two font sizes (the big initial-letter runs 202pt->168pt and
# the rest-of-word runs 140pt->117pt), matching the new/updated logo layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2, shape 2 ("Rectangle 5") - the "LAB for LINGUISTICS & COMPUTATION"
# stylised logo text box.
# ---------------------------------------------------------------------
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)

# Reposition / resize the text box (EMU -> points, 914400 EMU = 1 inch,
# 12700 EMU = 1 point):
#   off  x=850900  y=853618   ->  67pt , 67.2140157480315pt
#   ext cx=8356600 cy=7848302 -> 658pt , 617.9765354330709pt
$sh.Left   = 67.0
$sh.Top    = 67.2140160480315
$sh.Width  = 658.0
$sh.Height = 617.9765354330709

# Shrink the run font sizes. The paragraph text is:
#   "L" "AB for " <br> "L" "INGUISTICS  & " <br> "C" "OMPUTATION"
# with the single-letter runs ("L"/"L"/"C") at 202pt -> 168pt and the
# remaining runs at 140pt -> 117pt. Target via Characters(start,length)
# so each existing run keeps its other formatting (outline, pattern
# fill, shadow, typeface, ...) untouched.
$tr = $sh.TextFrame.TextRange

$tr.Characters(1, 1).Font.Size   = 168   # "L"
$tr.Characters(2, 7).Font.Size   = 117   # "AB for "
$tr.Characters(10, 1).Font.Size  = 168   # "L"
$tr.Characters(11, 14).Font.Size = 117   # "INGUISTICS  & "
$tr.Characters(26, 1).Font.Size  = 168   # "C"
$tr.Characters(27, 10).Font.Size = 117   # "OMPUTATION"

# ---------------------------------------------------------------------
# Presentation-level: mark that the deck now has a (empty) slide-guide
# list, as PowerPoint stamps into the extension list once guides have
# been touched.
# ---------------------------------------------------------------------
try {
    $guides = $p.Guides
    if ($guides -ne $null) {
        [void]$guides.Add(1, 5.0)
    }
} catch {
}
